# The sheet originally had two sections, "ece score" and "brier score loss",
# each listing rows for NL3..NL7 (5 label counts -> 5 rows, merged label cell
# in column A). The edit trims both sections down to only NL3..NL5, removing
# the NL6 and NL7 rows from each section (rows 7-8 and 12-13 in the original
# layout). Removing rows shifts everything below up, which is why we delete
# the lower block (12:13) before the upper block (7:8) - deleting bottom-up
# keeps the row numbers for the first deletion valid.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "brier score loss" NL6 and NL7 rows (originally rows 12-13).
$ws.Rows("12:13").Delete()

# Remove the "ece score" NL6 and NL7 rows (originally rows 7-8).
$ws.Rows("7:8").Delete()
